# Auto-generated edit script applying the 2022-12-04 YTD data update
# across 22 worksheets (129 cell value changes total).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E3").Value = 144
$ws.Range("B9").Value = 374
$ws.Range("C9").Value = 475
$ws.Range("D9").Value = 413
$ws.Range("E9").Value = 469
$ws.Range("F9").Value = 531
$ws.Range("I9").Value = 497
$ws.Range("B10").Value = 1328
$ws.Range("C10").Value = 1571
$ws.Range("D10").Value = 1787
$ws.Range("E10").Value = 2153
$ws.Range("F10").Value = 2086
$ws.Range("G10").Value = 888
$ws.Range("H10").Value = 597
$ws.Range("I10").Value = 845
$ws.Range("B11").Value = 1831
$ws.Range("C11").Value = 2205
$ws.Range("D11").Value = 2435
$ws.Range("E11").Value = 2849
$ws.Range("F11").Value = 2856
$ws.Range("G11").Value = 1553
$ws.Range("H11").Value = 1313
$ws.Range("I11").Value = 1679

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I7").Value = 9
$ws.Range("I9").Value = 26

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E7").Value = 53
$ws.Range("E9").Value = 154

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("C7").Value = 34
$ws.Range("C8").Value = 62
$ws.Range("F8").Value = 35
$ws.Range("C9").Value = 101
$ws.Range("F9").Value = 83

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E3").Value = 13
$ws.Range("F7").Value = 61
$ws.Range("B8").Value = 215
$ws.Range("C8").Value = 324
$ws.Range("D8").Value = 506
$ws.Range("E8").Value = 639
$ws.Range("F8").Value = 543
$ws.Range("H8").Value = 103
$ws.Range("B9").Value = 263
$ws.Range("C9").Value = 378
$ws.Range("D9").Value = 578
$ws.Range("E9").Value = 721
$ws.Range("F9").Value = 625
$ws.Range("H9").Value = 199

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("D6").Value = 12
$ws.Range("D7").Value = 19

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("D6").Value = 24
$ws.Range("D7").Value = 47
$ws.Range("D8").Value = 72

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("D5").Value = 19
$ws.Range("F7").Value = 24
$ws.Range("B8").Value = 75
$ws.Range("G8").Value = 90
$ws.Range("I21").Value = 26
$ws.Range("B28").Value = 100
$ws.Range("E32").Value = 154
$ws.Range("E35").Value = 22
$ws.Range("C36").Value = 101
$ws.Range("F36").Value = 83
$ws.Range("B53").Value = 263
$ws.Range("C53").Value = 378
$ws.Range("D53").Value = 578
$ws.Range("E53").Value = 721
$ws.Range("F53").Value = 625
$ws.Range("H53").Value = 199
$ws.Range("C62").Value = 26
$ws.Range("D65").Value = 72
$ws.Range("C68").Value = 21
$ws.Range("I74").Value = 43
$ws.Range("I76").Value = 46
$ws.Range("D77").Value = 55
$ws.Range("E77").Value = 73
$ws.Range("E78").Value = 51
$ws.Range("F78").Value = 53
$ws.Range("B86").Value = 14
$ws.Range("E92").Value = 36
$ws.Range("G92").Value = 30
$ws.Range("I92").Value = 38
$ws.Range("F95").Value = 60
$ws.Range("B97").Value = 32
$ws.Range("B99").Value = 1831
$ws.Range("C99").Value = 2205
$ws.Range("D99").Value = 2435
$ws.Range("E99").Value = 2849
$ws.Range("F99").Value = 2856
$ws.Range("G99").Value = 1553
$ws.Range("H99").Value = 1313
$ws.Range("I99").Value = 1679

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E5").Value = 43
$ws.Range("F5").Value = 42
$ws.Range("E6").Value = 51
$ws.Range("F6").Value = 53

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("B8").Value = 65
$ws.Range("B9").Value = 100

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I8").Value = 23
$ws.Range("I9").Value = 46

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("B6").Value = 11
$ws.Range("B7").Value = 14

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("E7").Value = 6
$ws.Range("G8").Value = 19
$ws.Range("I8").Value = 26
$ws.Range("E9").Value = 36
$ws.Range("G9").Value = 30
$ws.Range("I9").Value = 38

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("E6").Value = 16
$ws.Range("E7").Value = 22

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 16
$ws.Range("B7").Value = 32

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C5").Value = 23
$ws.Range("C6").Value = 26

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("D8").Value = 39
$ws.Range("E8").Value = 50
$ws.Range("D9").Value = 55
$ws.Range("E9").Value = 73

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("F5").Value = 53
$ws.Range("F6").Value = 60

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("F5").Value = 7
$ws.Range("F7").Value = 24

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("C6").Value = 19
$ws.Range("C7").Value = 21

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("B6").Value = 22
$ws.Range("G7").Value = 56
$ws.Range("B8").Value = 75
$ws.Range("G8").Value = 90
